# Database refresh for the "Overview" income-statement sheet.
#
# The sheet keeps a rolling 5-period window (columns D:H) of yearly figures.
# A new fiscal year (1401/12) has been published, so the window slides one
# column to the left (the oldest period, 1396/12, falls off the front) and
# the newly opened rightmost column (H) is populated with the freshly
# published period's figures. This also changes how the per-period "price"
# figures are derived (read_price algorithm), which is why more than a
# simple shift shows up in a few rows (H26 "سرمایه", and the recomputed
# H-column totals throughout).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "E", "F", "G", "H")

function Set-RowValues([int]$row, [object[]]$values) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $values[$i]
    }
}

# Row 8 - "دوره مالی" (period labels)
Set-RowValues 8 @(
    "12 ماهه منتهی به 1397/12",
    "12 ماهه منتهی به 1398/12",
    "12 ماهه منتهی به 1399/12",
    "12 ماهه منتهی به 1400/12",
    "12 ماهه منتهی به 1401/12"
)

# Row 9 - "تاریخ انتشار" (publish dates)
Set-RowValues 9 @(
    "1399-02-11 (8)",
    "1400-02-08 (8)",
    "1401-02-06 (9)",
    "1402-02-10 (8)",
    "1402-02-10 (2)"
)

# Row 11 - "فروش" (sales)
Set-RowValues 11 @(1717908, 2641417, 3409060, 8310129, 11031555)

# Row 12 - "بهای تمام شده کالای فروش رفته" (cost of goods sold)
Set-RowValues 12 @(-1370356, -2110736, -2484486, -6576671, -8781736)

# Row 13 - "سود (زیان) ناخالص" (gross profit)
Set-RowValues 13 @(0, 530681, 924574, 1733458, 2249819)

# Row 14 - "هزینه های عمومی, اداری و تشکیلاتی" (G&A expenses)
Set-RowValues 14 @(-56056, -67946, -117229, -163847, -215884)

# Row 15 - "هزینه کاهش ارزش دریافتنی‌‏ها (هزینه استثنایی)" (impairment expense)
Set-RowValues 15 @(-17391, 0, 0, 0, 0)

# Row 16 - "خالص سایر درامدها (هزینه ها) ی عملیاتی" (other operating income/expense, net)
Set-RowValues 16 @(-5549, 5352, 30553, 115502, 33125)

# Row 17 - "سود (زیان) عملیاتی" (operating profit)
Set-RowValues 17 @(0, 468087, 837898, 1685113, 2067060)

# Row 18 - "هزینه های مالی" (finance costs)
Set-RowValues 18 @(-41425, -16531, -9111, -14203, -13393)

# Row 19 - "خالص سایر درامدها و هزینه های غیرعملیاتی" (other non-operating income/expense, net)
Set-RowValues 19 @(65246, -347, -44015, 17174, -12868)

# Row 20 - "سود (زیان) خالص عملیات در حال تداوم قبل از مالیات" (pre-tax profit from continuing ops)
Set-RowValues 20 @(0, 451209, 784772, 1688084, 2040799)

# Row 21 - "مالیات" (tax)
Set-RowValues 21 @(-69976, -101089, -113671, -282485, -209718)

# Row 22 - "سود (زیان) خالص عملیات در حال تداوم" (net profit from continuing ops)
Set-RowValues 22 @(0, 350120, 671101, 1405599, 1831081)

# Row 23 - "سود (زیان) عملیات متوقف شده پس از اثر مالیاتی" (discontinued ops, after tax)
Set-RowValues 23 @(0, 0, 0, 0, 0)

# Row 24 - "سود (زیان) خالص" (net profit)
Set-RowValues 24 @(0, 350120, 671101, 1405599, 1831081)

# Row 25 - "سود هر سهم پس از کسر مالیات" (EPS after tax)
Set-RowValues 25 @(0, 738, 1414, 2962, 1927)

# Row 26 - "سرمایه" (capital) - capital increase lands in the new column
Set-RowValues 26 @(474522, 474522, 474522, 474522, 950000)

# Row 27 - "سود هر سهم بر اساس آخرین سرمایه" (EPS based on latest capital)
Set-RowValues 27 @(0, 369, 706, 1480, 1927)
